$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E (상태/Status) and Column F (진척도/Progress) for the WBS rows ---
# 완료 (Done) rows
$ws.Range("E13").Value = "완료"
$ws.Range("F13").Value = 1
$ws.Range("E14").Value = "완료"
$ws.Range("F14").Value = 1
$ws.Range("E15").Value = "완료"
$ws.Range("F15").Value = 1
$ws.Range("E16").Value = "완료"
$ws.Range("F16").Value = 1

# 진행중 (In progress) rows
$ws.Range("E17").Value = "진행중"
$ws.Range("F17").Formula = "=(F$19+F$20)/2"
$ws.Range("E18").Value = "진행중"
$ws.Range("F18").Formula = "=(F$19+F$20)/2"
$ws.Range("E19").Value = "진행중"
$ws.Range("F19").Value = 0.4
$ws.Range("E20").Value = "진행중"
$ws.Range("F20").Value = 0.3
$ws.Range("E21").Value = "진행중"
$ws.Range("F21").Formula = "=(F23+F24+F25+F26+F28+F29+F31+F32)/8"
$ws.Range("E22").Value = "진행중"
$ws.Range("F22").Formula = "=(F23+F24+F25+F26)/4"
$ws.Range("E23").Value = "진행중"
$ws.Range("F23").Value = 0
$ws.Range("E24").Value = "진행중"
$ws.Range("F24").Value = 0.7
$ws.Range("E25").Value = "진행중"
$ws.Range("F25").Value = 0.8
$ws.Range("E26").Value = "진행중"
$ws.Range("F26").Value = 0.5
$ws.Range("E27").Value = "진행중"
$ws.Range("F27").Formula = "=(F28+F29)/2"
$ws.Range("E28").Value = "진행중"
$ws.Range("F28").Value = 0
$ws.Range("E29").Value = "진행중"
$ws.Range("F29").Value = 0.3
$ws.Range("E30").Value = "진행중"
$ws.Range("F30").Formula = "=(F31+F32)/2"
$ws.Range("E31").Value = "진행중"
$ws.Range("F31").Value = 0
$ws.Range("E32").Value = "진행중"
$ws.Range("F32").Value = 0.5
$ws.Range("E33").Value = "진행중"
$ws.Range("F33").Formula = "=F35"
$ws.Range("E34").Value = "진행중"
$ws.Range("F34").Value = 0.3
$ws.Range("E35").Value = "진행중"
$ws.Range("F35").Value = 0.3

# 예정 (Planned) rows
$ws.Range("E36").Value = "예정"
$ws.Range("F36").Value = 0
$ws.Range("E37").Value = "예정"
$ws.Range("F37").Value = 0
$ws.Range("E38").Value = "예정"
$ws.Range("F38").Value = 0
$ws.Range("E39").Value = "예정"
$ws.Range("F39").Value = 0
$ws.Range("E40").Value = "예정"
$ws.Range("F40").Value = 0
$ws.Range("E41").Value = "예정"
$ws.Range("F41").Value = 0
$ws.Range("E42").Value = "예정"
$ws.Range("F42").Value = 0
$ws.Range("E43").Value = "예정"
$ws.Range("F43").Value = 0

# Recalculate so formula cells (F17, F18, F21, F22, F27, F30, F33) carry fresh cached values
$excel.Calculate()

# --- Restore the author's last on-screen selection ---
[void]$ws.Range("F33").Select()
